$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the "Action recommandée" (column E) value from "Oui" to "Non"
# for the rows affected by the edit.
$rows = @(3, 4, 5, 6, 7, 21, 23, 24)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "Non"
}

# Update the active selection on the sheet from E32 to E5.
$ws.Range("E5").Select()
